$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so values such as
# "236.29" or "1.000" are not silently reinterpreted as numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.399.58'
$ws.Range("E2").Value = '  -1.93%  '

$ws.Range("D3").Value = '1.664.48'
$ws.Range("E3").Value = '  -3.83%  '

$ws.Range("D4").Value = '0.9977'
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = '236.29'
$ws.Range("E5").Value = '  -3.69%  '

$ws.Range("D6").Value = '0.9988'
$ws.Range("E6").Value = '  -0.16%  '

$ws.Range("D7").Value = '0.4806'
$ws.Range("E7").Value = '  -4.37%  '

$ws.Range("D8").Value = '0.2605'
$ws.Range("E8").Value = '  -4.12%  '

$ws.Range("D9").Value = '0.06146'
$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("D10").Value = '0.07069'
$ws.Range("E10").Value = '  -2.30%  '

$ws.Range("D11").Value = '1.656.20'
$ws.Range("E11").Value = '  -4.52%  '

$ws.Range("D12").Value = '14.71'
$ws.Range("E12").Value = '  -2.93%  '

$ws.Range("D13").Value = '0.5890'
$ws.Range("E13").Value = '  -9.76%  '

$ws.Range("D14").Value = '4.372'
$ws.Range("E14").Value = '  -8.17%  '

$ws.Range("D15").Value = '74.35'
$ws.Range("E15").Value = '  -3.43%  '

$ws.Range("E16").Value = '  -0.10%  '

$ws.Range("D17").Value = '0.9985'
$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("D18").Value = '25.383.32'
$ws.Range("E18").Value = '  -2.02%  '

$ws.Range("D19").Value = '0.000006698'
$ws.Range("E19").Value = '  -1.93%  '

$ws.Range("E20").Value = '  -4.16%  '

$ws.Range("D21").Value = '1.869.85'
$ws.Range("E21").Value = '  -4.57%  '

$ws.Range("D22").Value = '4.400'
$ws.Range("E22").Value = '  -4.16%  '

$ws.Range("D23").Value = '8.616'
$ws.Range("E23").Value = '  -2.20%  '

$ws.Range("D24").Value = '5.320'
$ws.Range("E24").Value = '  -2.89%  '

$ws.Range("D25").Value = '133.97'
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("E26").Value = '  -0.95%  '

$ws.Range("E27").Value = '  -2.15%  '

$ws.Range("D28").Value = '104.71'
$ws.Range("E28").Value = '  -0.40%  '

$ws.Range("D29").Value = '1.685'
$ws.Range("E29").Value = '  -5.81%  '

$ws.Range("D30").Value = '3.994'
$ws.Range("E30").Value = '  +0.61%  '

$ws.Range("D31").Value = '3.612'
$ws.Range("E31").Value = '  -2.39%  '

$ws.Range("D32").Value = '0.07649'
$ws.Range("E32").Value = '  -5.92%  '

$ws.Range("D33").Value = '0.04370'
$ws.Range("E33").Value = '  -7.63%  '

$ws.Range("D34").Value = '0.9981'
$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("E35").Value = '  -2.06%  '

$ws.Range("D36").Value = '0.6032'
$ws.Range("E36").Value = '  -1.74%  '

$ws.Range("D37").Value = '0.9423'
$ws.Range("E37").Value = '  -5.35%  '

$ws.Range("D38").Value = '2.619'
$ws.Range("E38").Value = '  -4.50%  '

$ws.Range("D39").Value = '0.8530'
$ws.Range("E39").Value = '  -3.72%  '

$ws.Range("D40").Value = '0.9994'
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("B41").Value = 'PaxosStandard'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D41").Value = '0.9986'
$ws.Range("E41").Value = '  -0.28%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.01501'
$ws.Range("E42").Value = '  -6.82%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '1.819'
$ws.Range("E43").Value = '  -7.28%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '98.60'
$ws.Range("E44").Value = '  -3.03%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3748'
$ws.Range("E45").Value = '  -4.06%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '4.693'
$ws.Range("E46").Value = '  -6.38%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1113'
$ws.Range("E47").Value = '  -5.78%  '

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '6.198'
$ws.Range("E48").Value = '  -2.71%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.05245'
$ws.Range("E49").Value = '  -0.64%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '29.42'
$ws.Range("E50").Value = '  -4.61%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.218'
$ws.Range("E51").Value = '  -1.75%  '

# Restore the original (default) style on the Price column now that
# the text values have been written, so no stray style index remains.
$priceRange.Style = "Normal"
